# Fill in the "IkariNote03 (windows64)" results (column B) on the
# "Ejecucion mas rápida" sheet - data captured from "windows note03".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ejecucion mas rápida")

$ws.Range("B6").Formula  = "=8 *376318.99"
$ws.Range("B7").Formula  = "=4 * 585297.42"
$ws.Range("B8").Formula  = "= 4 * 268148.63"
$ws.Range("B9").Formula  = "=2*220681.42"
$ws.Range("B10").Value   = 588878.35
$ws.Range("B11").Value   = 594142.91
$ws.Range("B12").Value   = 200044.69
$ws.Range("B13").Value   = 203012.79
$ws.Range("B14").Value   = 455707.54
$ws.Range("B15").Value   = 126068.81
$ws.Range("B16").Value   = 119128.89

# Update the saved selection to match the author's final cursor position.
$ws.Range("B14").Select()
